$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab/label from "BrassA-HW15.xpc" to "BrassA"
$ws.Name = "BrassA"

# Append a new row of averaged intensity data (row 16), matching the
# HKL label/style pattern used by row 15 (HexGrid-60degTilt5degRes).
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = $ws.Range("B15").Value2

$ws.Range("C16").Value2 = 1.381176291822647
$ws.Range("D16").Value2 = 1.012336909755702
$ws.Range("E16").Value2 = 0.9419819153733993
$ws.Range("F16").Value2 = 0.9383915040732334
$ws.Range("G16").Value2 = 1.381176291822647
$ws.Range("H16").Value2 = 1.012336909755702
$ws.Range("I16").Value2 = 1.019680040198932
$ws.Range("J16").Value2 = 0.8337931478668805
$ws.Range("K16").Value2 = 1.075950142525188
$ws.Range("L16").Value2 = 0.934283842502692
$ws.Range("M16").Value2 = 1.381176291822647
$ws.Range("N16").Value2 = 0.9771594125645509
$ws.Range("O16").Value2 = 1.068471655256246
$ws.Range("P16").Value2 = 1.017199224264834
